$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33 (item id 5512)
$ws.Range("H33").Value = 218.45454
$ws.Range("I33").Value = 218.45454
$ws.Range("K33").Value = 218.45454
$ws.Range("M33").Value = 10.54545999999999

# Row 40 (item id 5505)
$ws.Range("H40").Value = 1677.6
$ws.Range("I40").Value = 997
$ws.Range("K40").Value = 997
$ws.Range("M40").Value = -822

# Row 43 (item id 5472)
$ws.Range("H43").Value = 1490
$ws.Range("I43").Value = 983.3333
$ws.Range("J43").Value = 2250
$ws.Range("K43").Value = 983.3333
$ws.Range("L43").Value = 2250
$ws.Range("M43").Value = -914.3333
$ws.Range("N43").Value = -2388

# Row 51 (item id 5486)
$ws.Range("H51").Value = 3084.8044
$ws.Range("I51").Value = 2864.3096
$ws.Range("K51").Value = 2864.3096
$ws.Range("M51").Value = -2380.3096

# Row 55 (item id 5517)
$ws.Range("H55").Value = 143.66667
$ws.Range("I55").Value = 147.875
$ws.Range("J55").Value = 135.25
$ws.Range("K55").Value = 147.875
$ws.Range("L55").Value = 135.25
$ws.Range("M55").Value = 66.125
$ws.Range("N55").Value = -563.25

# Row 131 (item id 36108)
$ws.Range("H131").Value = 2362.375
$ws.Range("I131").Value = 2485.5715
$ws.Range("J131").Value = 1500
$ws.Range("K131").Value = 7456.7145
$ws.Range("L131").Value = 4500
$ws.Range("M131").Value = -2416.7145
$ws.Range("N131").Value = -14580

# Row 132 (item id 44049)
$ws.Range("H132").Value = 1488.2354
$ws.Range("I132").Value = 1488.2354
$ws.Range("K132").Value = 4464.706200000001
$ws.Range("M132").Value = -1934.706200000001

# Row 137 (item id 44013)
$ws.Range("H137").Value = 8512.134
$ws.Range("I137").Value = 2098.2
$ws.Range("J137").Value = 21340
$ws.Range("K137").Value = 6294.599999999999
$ws.Range("L137").Value = 64020
$ws.Range("M137").Value = -3744.599999999999
$ws.Range("N137").Value = -69120

# Row 138 (item id 44169)
$ws.Range("H138").Value = 5470
$ws.Range("I138").Value = 6629.769
$ws.Range("J138").Value = 4676.4736
$ws.Range("K138").Value = 19889.307
$ws.Range("L138").Value = 14029.4208
$ws.Range("M138").Value = -14749.307
$ws.Range("N138").Value = -24309.4208


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (item id 44147)
$ws.Range("H32").Value = 6272.3
$ws.Range("I32").Value = 4076.9363
$ws.Range("K32").Value = 4076.9363
$ws.Range("M32").Value = -3789.9363

# Row 45 (item id 27714)
$ws.Range("H45").Value = 4986
$ws.Range("I45").Value = 3300
$ws.Range("J45").Value = 5226.857
$ws.Range("K45").Value = 3300
$ws.Range("L45").Value = 5226.857
$ws.Range("M45").Value = -2923
$ws.Range("N45").Value = -5980.857

# Row 74 (item id 44000)
$ws.Range("H74").Value = 7933.3335
$ws.Range("I74").Value = 3535.6086
$ws.Range("J74").Value = 14255.0625
$ws.Range("K74").Value = 3535.6086
$ws.Range("L74").Value = 14255.0625
$ws.Range("M74").Value = -2661.6086
$ws.Range("N74").Value = -16003.0625

# Row 77 (item id 44000)
$ws.Range("H77").Value = 7933.3335
$ws.Range("I77").Value = 3535.6086
$ws.Range("J77").Value = 14255.0625
$ws.Range("K77").Value = 17678.043
$ws.Range("L77").Value = 71275.3125
$ws.Range("M77").Value = -13310.043
$ws.Range("N77").Value = -80011.3125

# Row 122 (item id 36168)
$ws.Range("H122").Value = 2315.4285
$ws.Range("I122").Value = 2179.2727
$ws.Range("J122").Value = 2814.6667
$ws.Range("K122").Value = 6537.8181
$ws.Range("L122").Value = 8444.000100000001
$ws.Range("M122").Value = -4087.8181
$ws.Range("N122").Value = -13344.0001

# Row 124 (item id 34252)
$ws.Range("H124").Value = 81000
$ws.Range("J124").Value = 81000
$ws.Range("L124").Value = 81000
$ws.Range("N124").Value = -90820

# Row 125 (item id 34251)
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("N125").ClearContents()


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107 (item id 27706)
$ws.Range("H107").Value = 1829.875
$ws.Range("I107").Value = 1856.5
$ws.Range("K107").Value = 1856.5
$ws.Range("M107").Value = 63.5

# Row 112 (item id 25788)
$ws.Range("H112").Value = 41234.5
$ws.Range("J112").Value = 41234.5
$ws.Range("L112").Value = 41234.5
$ws.Range("N112").Value = -44188.5

# Row 134 (item id 43998)
$ws.Range("H134").Value = 7148.606
$ws.Range("I134").Value = 3901.8845
$ws.Range("K134").Value = 11705.6535
$ws.Range("M134").Value = -9170.6535


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22 (item id 5367)
$ws.Range("H22").Value = 456606.34
$ws.Range("I22").Value = 656213.5600000001
$ws.Range("J22").Value = 2953.5454
$ws.Range("K22").Value = 656213.5600000001
$ws.Range("L22").Value = 2953.5454
$ws.Range("M22").Value = -655863.5600000001
$ws.Range("N22").Value = -3653.5454

# Row 31 (item id 44023)
$ws.Range("H31").Value = 2758.4666
$ws.Range("I31").Value = 2347.75
$ws.Range("J31").Value = 3227.8572
$ws.Range("K31").Value = 2347.75
$ws.Range("L31").Value = 3227.8572
$ws.Range("M31").Value = -2052.75
$ws.Range("N31").Value = -3817.8572

# Row 34 (item id 44023)
$ws.Range("H34").Value = 2758.4666
$ws.Range("I34").Value = 2347.75
$ws.Range("J34").Value = 3227.8572
$ws.Range("K34").Value = 2347.75
$ws.Range("L34").Value = 3227.8572
$ws.Range("M34").Value = -2145.75
$ws.Range("N34").Value = -3631.8572

# Row 99 (item id 36198)
$ws.Range("H99").Value = 16884.2
$ws.Range("J99").Value = 3341.5557
$ws.Range("L99").Value = 3341.5557
$ws.Range("N99").Value = -6337.5557

# Row 126 (item id 36198)
$ws.Range("H126").Value = 16884.2
$ws.Range("J126").Value = 3341.5557
$ws.Range("L126").Value = 10024.6671
$ws.Range("N126").Value = -14964.6671

# Row 129 (item id 35378)
$ws.Range("H129").Value = 60000
$ws.Range("J129").Value = 60000
$ws.Range("L129").Value = 60000
$ws.Range("N129").Value = -70000

# Row 132 (item id 44019)
$ws.Range("H132").Value = 4182.375
$ws.Range("I132").Value = 3972.7144
$ws.Range("K132").Value = 11918.1432
$ws.Range("M132").Value = -9388.143199999999


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2 (item id 4847)
$ws.Range("H2").Value = 119.6
$ws.Range("I2").Value = 84.5
$ws.Range("K2").Value = 507
$ws.Range("M2").Value = -394

# Row 13 (item id 4657)
$ws.Range("H13").Value = 113.5
$ws.Range("I13").Value = 18.5
$ws.Range("J13").Value = 208.5
$ws.Range("K13").Value = 55.5
$ws.Range("L13").Value = 625.5
$ws.Range("M13").Value = 112.5
$ws.Range("N13").Value = -961.5

# Row 15 (item id 4661)
$ws.Range("H15").Value = 80
$ws.Range("I15").Value = 34.6
$ws.Range("J15").Value = 125.4
$ws.Range("K15").Value = 103.8
$ws.Range("L15").Value = 376.2
$ws.Range("M15").Value = 36.19999999999999
$ws.Range("N15").Value = -656.2

# Row 34 (item id 4749)
$ws.Range("H34").Value = 2375.3333
$ws.Range("I34").Value = 850.6667
$ws.Range("J34").Value = 3900
$ws.Range("K34").Value = 2552.0001
$ws.Range("L34").Value = 11700
$ws.Range("M34").Value = -2468.0001
$ws.Range("N34").Value = -11868

# Row 69 (item id 12850)
$ws.Range("H69").Value = 2584.5
$ws.Range("I69").Value = 2500
$ws.Range("J69").Value = 2669
$ws.Range("K69").Value = 7500
$ws.Range("L69").Value = 8007
$ws.Range("M69").Value = -6689
$ws.Range("N69").Value = -9629

# Row 72 (item id 12850)
$ws.Range("H72").Value = 2584.5
$ws.Range("I72").Value = 2500
$ws.Range("J72").Value = 2669
$ws.Range("K72").Value = 22500
$ws.Range("L72").Value = 24021
$ws.Range("M72").Value = -18444
$ws.Range("N72").Value = -32133

# Row 87 (item id 12864)
$ws.Range("H87").Value = 9798.200000000001
$ws.Range("I87").Value = 9798.200000000001
$ws.Range("K87").Value = 29394.6
$ws.Range("M87").Value = -28146.6

# Row 90 (item id 12864)
$ws.Range("H90").Value = 9798.200000000001
$ws.Range("I90").Value = 9798.200000000001
$ws.Range("K90").Value = 88183.8
$ws.Range("M90").Value = -81943.8

# Row 134 (item id 44074)
$ws.Range("H134").Value = 1588.5714
$ws.Range("I134").Value = 1588.5714
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4765.7142
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 68 (item id 10659)
$ws.Range("H68").Value = 35000
$ws.Range("J68").Value = 35000
$ws.Range("L68").Value = 35000
$ws.Range("N68").Value = -36622

# Row 71 (item id 10659)
$ws.Range("H71").Value = 35000
$ws.Range("J71").Value = 35000
$ws.Range("L71").Value = 105000
$ws.Range("N71").Value = -113112

# Row 75 (item id 11008)
$ws.Range("H75").Value = 46131
$ws.Range("J75").Value = 46131
$ws.Range("L75").Value = 46131
$ws.Range("N75").Value = -47879

# Row 78 (item id 11008)
$ws.Range("H78").Value = 46131
$ws.Range("J78").Value = 46131
$ws.Range("L78").Value = 138393
$ws.Range("N78").Value = -147129

# Row 123 (item id 34150)
$ws.Range("H123").Value = 64925.5
$ws.Range("J123").Value = 64925.5
$ws.Range("L123").Value = 64925.5
$ws.Range("N123").Value = -69825.5

# Row 126 (item id 36184)
$ws.Range("H126").Value = 3378.25
$ws.Range("I126").Value = 3173
$ws.Range("J126").Value = 3501.4
$ws.Range("K126").Value = 9519
$ws.Range("L126").Value = 10504.2
$ws.Range("M126").Value = -7049
$ws.Range("N126").Value = -15444.2

# Row 132 (item id 44008)
$ws.Range("H132").Value = 8538.857
$ws.Range("I132").Value = 9780.412
$ws.Range("J132").Value = 3262.25
$ws.Range("K132").Value = 29341.236
$ws.Range("L132").Value = 9786.75
$ws.Range("M132").Value = -26811.236
$ws.Range("N132").Value = -14846.75


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22 (item id 5277)
$ws.Range("H22").Value = 3082.6296
$ws.Range("I22").Value = 1449.2858
$ws.Range("J22").Value = 3654.3
$ws.Range("K22").Value = 1449.2858
$ws.Range("L22").Value = 3654.3
$ws.Range("M22").Value = -1154.2858
$ws.Range("N22").Value = -4244.3

# Row 27 (item id 5277)
$ws.Range("H27").Value = 3082.6296
$ws.Range("I27").Value = 1449.2858
$ws.Range("J27").Value = 3654.3
$ws.Range("K27").Value = 1449.2858
$ws.Range("L27").Value = 3654.3
$ws.Range("M27").Value = -1342.2858
$ws.Range("N27").Value = -3868.3

# Row 46 (item id 5282)
$ws.Range("H46").Value = 3971.52
$ws.Range("I46").Value = 325.5
$ws.Range("K46").Value = 325.5
$ws.Range("M46").Value = -137.5

# Row 55 (item id 5284)
$ws.Range("H55").Value = 1387
$ws.Range("I55").Value = 1327.5714
$ws.Range("J55").Value = 1426.619
$ws.Range("K55").Value = 1327.5714
$ws.Range("L55").Value = 1426.619
$ws.Range("M55").Value = -1154.5714
$ws.Range("N55").Value = -1772.619

# Row 111 (item id 25820)
$ws.Range("H111").Value = 48387
$ws.Range("J111").Value = 48387
$ws.Range("L111").Value = 48387
$ws.Range("N111").Value = -56567

# Row 122 (item id 36247)
$ws.Range("H122").Value = 7832.3335
$ws.Range("I122").Value = 4994
$ws.Range("J122").Value = 8400
$ws.Range("K122").Value = 14982
$ws.Range("L122").Value = 25200
$ws.Range("M122").Value = -12532
$ws.Range("N122").Value = -30100

# Row 132 (item id 44058)
$ws.Range("H132").Value = 4590.2
$ws.Range("I132").Value = 4590.2
$ws.Range("K132").Value = 13770.6
$ws.Range("M132").Value = -11240.6


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 40 (item id 3601)
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

